# Add the 2025 national-contest ("国赛") row to the classification sheet,
# mirroring the other yearly rows already present (row 22 -> year 2025).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Match the row height / banding used by the other data rows (2-21).
$ws.Rows.Item(22).RowHeight = 30

# Year column.
$ws.Range("A22").Value = 2025

# Write the new text cells in the same order the original author typed
# them in (this governs the order new entries are appended to the shared
# string table, matching the upstream workbook byte-for-byte in content).
$ws.Range("B22").Value = "D_简易以太网双绞线测试仪"
$ws.Range("I22").Value = "F_简易自动接收机"
$ws.Range("D22").Value = "A_能量回馈的变流器负载试验装置`nB_单相有源电力滤波实验装置"
$ws.Range("H22").Value = "H_野生动物巡查系统"
$ws.Range("F22").Value = "G_电路模型探究装置"
$ws.Range("G22").Value = "E_简易自行瞄准装置`nC_基于单目视觉的目标物测量装置"

# D22/G22 hold multi-line text, same as the analogous multi-entry cells in
# earlier rows, so wrap them like those.
$ws.Range("D22").WrapText = $true
$ws.Range("G22").WrapText = $true

# Leave the selection on the newly added row, like the source workbook.
$ws.Range("G19").Select()
